$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I (I0) and J (IF), matching the style of the
# existing header cells (bold, thin border, centered/top-aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160

# Data values for columns I (I0) and J (IF), rows 2-34
$data = @(
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(11, 11),
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(7, 8),
    @(9, 10),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(9, 9),
    @(6, 7),
    @(1, 3),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(1, 2),
    @(9, 9),
    @(6, 6),
    @(4, 5),
    @(7, 7),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
